$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")
$ws.Cells.Item(27, 5).Value = ""
